$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.841.25"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.350.21"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'548.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'132.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "'0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").Value = "'5.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.24%  "
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "'0.352"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "'23.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "2.766.77"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "57.786.95"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "2.330.12"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "'11.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.68%  "
$ws.Range("D19").Value = "'4.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'328.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'6.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'63.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'0.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("E27").Value = "  -5.29%  "
$ws.Range("D28").Value = "'1.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'171.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "'6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'18.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'4.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.424"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.22%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'141.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'289.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'0.0951"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "Polygon"
$ws.Range("C44").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D44").Value = "'0.418"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.55%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0513"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'18.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").Value = "'0.0221"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "'11.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'4.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'0.941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
